$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column K: "Punteggio medio" (average score) ---

# Give K1:K2 the same look as the "Domanda" header cell (A2): copy its
# formatting (font/fill/alignment) onto the new header cells before
# merging, then set the label text.
$ws.Range("A2").Copy()
$ws.Range("K1:K2").PasteSpecial(-4122)
$ws.Range("K1").Value2 = "Punteggio medio"
$ws.Range("K1:K2").Merge()

# Average of the four "Voto totale" cells
$ws.Range("K13").Formula = "=AVERAGE(D13,F13,H13,J13)"

# Widen the new column similarly to the existing score columns
$ws.Columns.Item(11).ColumnWidth = 16.5

# Leave the sheet scrolled down with the new total selected
$ws.Range("K13").Select()
